{"js": "// The _GoBack bookmark that Word drops at the last edit position was\n// previously sitting alone in a trailing empty paragraph at the very end\n// of the document. Word's \"last cursor position\" tracking moved it back\n// to the start of the second paragraph (the one beginning \"Dalam project\n// ...\") and the now-empty trailing paragraph was removed.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Insert the _GoBack bookmark at the very start of the second paragraph\n// (\"Dalam project ...\"), right before its first run.\nconst target = items[1];\ntarget.getRange(\"Start\").insertBookmark(\"_GoBack\");\n\n// Remove the trailing empty paragraph (it used to hold the _GoBack\n// bookmark by itself).\nconst last = items[items.length - 1];\nlast.delete();\n\nawait context.sync();\n", "ps1": "# The _GoBack bookmark that Word drops at the last edit position used to\n# sit by itself in a trailing empty paragraph at the very end of the\n# document. Move it back to the start of the second paragraph (the one\n# that begins \"Dalam project ...\") and remove the now-empty trailing\n# paragraph.\n\n$d = $word.ActiveDocument\n\n# Remove the trailing empty paragraph that currently holds the old\n# _GoBack bookmark (deleting the paragraph also removes that bookmark).\n$last = $d.Paragraphs.Last\n$last.Range.Delete()\n\n# Collapse a range to the very start of the second paragraph and drop a\n# fresh _GoBack bookmark there.\n$target = $d.Paragraphs(2)\n$r = $target.Range\n$r.Collapse(1)\n$d.Bookmarks.Add(\"_GoBack\", $r)\n"}
